# Apply the checklist update: insert a new "WSTG-CLNT-15" row into the
# "Testing Checklist" sheet right before the existing "API Testing" section
# (i.e. at row 133), pushing the API Testing section down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Checklist")

# 1) Insert a new row at 133. This shifts the old rows 133-138 down to 134-139
#    (the "API Testing" header + WSTG-APIT-01/02/99 rows + trailing blank
#    separator row all move down by one row, unchanged).
$ws.Rows.Item(133).Insert()

# 2) Copy the formatting of the WSTG-APIT-01 content row (now at row 136,
#    after the shift) onto the freshly inserted row 133, since it uses the
#    same visual style (borders/fonts/alignment) as the new WSTG-CLNT-15
#    entry we are about to add.
$ws.Range("A136:F136").Copy()
$ws.Range("A133:F133").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Set the row height to match the new content (99pt, matching the other
#    multi-line WSTG content rows).
$ws.Rows.Item(133).RowHeight = 99

# 4) Populate the new row's content.
$ws.Cells.Item(133, 1).Value = ""
$ws.Cells.Item(133, 2).Value = "WSTG-CLNT-15"
$ws.Cells.Item(133, 3).Formula = '=HYPERLINK("https://owasp.org/www-project-web-security-testing-guide/latest/4-Web_Application_Security_Testing/11-Client-side_Testing/15-Testing_for_Client-Side_Template_Injection", "Testing for Client-side Template Injection")'

$objectives = "- Identify the client-side framework and its version used by the application." + "`n" + "- Detect injection points where user input is reflected into the DOM and processed by the template engine." + "`n" + "- Assess if the injection allows for arbitrary JavaScript execution (XSS) via the template syntax."
$ws.Cells.Item(133, 4).Value = $objectives

$ws.Cells.Item(133, 5).Value = "Not Started"
$ws.Cells.Item(133, 6).Value = ""

# 5) The conditional formatting that covered B4:F138 needs to grow to
#    B4:F139 to keep covering the whole table (it does not auto-expand
#    when a row is inserted at the very edge of its range).
$newCfRange = $ws.Range("B4:F139")
$conditions = $ws.Cells.FormatConditions
for ($i = 1; $i -le $conditions.Count; $i++) {
    $cond = $conditions.Item($i)
    $addr = $cond.AppliesTo.Address()
    if ($addr -eq "`$B`$4:`$F`$138") {
        $cond.ModifyAppliesToRange($newCfRange)
    }
}

# 6) Data validation (the "Status" dropdown) should also cover the new
#    E133 cell, in addition to the cells that already had it (whose
#    references were shifted automatically by the row insert above).
$ws.Range("E133").Validation.Add(3, 1, 1, "Not Started,Pass,Issues,N/A")
